$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
Write-Output $ws.Name
